# Weekly update: a new week of "Espinaca" price data for Feria Lagunitas de
# Puerto Montt is inserted right before the current row 45, pushing all the
# following historical rows (old 45..57) down by one (new 46..58).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 45..57 down to 46..58 and create a fresh, blank row 45.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with this week's data point.
$ws.Cells.Item(45, 1).Value = 4
$ws.Cells.Item(45, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(45, 3).Value = "Los Lagos"
$ws.Cells.Item(45, 4).Value = 45009
$ws.Cells.Item(45, 5).Value = 10
$ws.Cells.Item(45, 6).Value = 100112012
$ws.Cells.Item(45, 7).Value = "Espinaca"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 25
$ws.Cells.Item(45, 11).Value = 15000
$ws.Cells.Item(45, 12).Value = 15000
$ws.Cells.Item(45, 13).Value = 15000
$ws.Cells.Item(45, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(45, 15).Value = "Región Metropolitana"
$ws.Cells.Item(45, 16).Value = 1500
$ws.Cells.Item(45, 17).Value = 10
$ws.Cells.Item(45, 18).Value = "Hortaliza"

# Match the date format used by the rest of column D.
$ws.Cells.Item(45, 4).NumberFormat = $ws.Cells.Item(46, 4).NumberFormat()
